# EPICP-1 added ident to DD_EPICP_TRACY
$wb = $excel.ActiveWorkbook

$wsVars = $wb.Worksheets.Item("Variables")
$wsCats = $wb.Worksheets.Item("Categories")

# Insert a new row at row 2 on the "Variables" sheet, pushing existing rows down.
$wsVars.Rows.Item(2).Insert()

# Copy header-row formatting (row 1) onto the newly inserted row 2 so it
# matches the bold/centered "customFormat" look used for the header.
$wsVars.Rows.Item(1).Copy()
$wsVars.Rows.Item(2).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$wsVars.Application.CutCopyMode = $false

# Fill in the new variable: ident / participant identifier / integer
$wsVars.Range("A2").Value = ""
$wsVars.Range("B2").Value = "ident"
$wsVars.Range("C2").Value = "participant identifier"
$wsVars.Range("D2").Value = "integer"

# Restore default (non-wrapped) row height on what is now row 29 (previously
# row 28, which had an explicit ht="30" due to wrapped text).
$wsVars.Rows.Item(29).RowHeight = $wsVars.Rows.Item(3).RowHeight

# Update the saved selections/view state to match the edited workbook.
$wsVars.Range("H9").Select()

$wsCats.Activate()
$wsCats.Range("A2:XFD2").Select()
$wsCats.Application.ActiveWindow.ScrollRow = 19
